$d = $word.ActiveDocument

# ------------------------------------------------------------------------
# Target change: the "{{Telefono}}" placeholder becomes "{telefono}".
#
# In the canonical OOXML the placeholder was spread across six runs:
#   "{{" | "Tel" | "e" | "fono" | "}" | "}"
# and must end up, after the edit, spread across six (different) runs:
#   "{" | "t" | "el" | "e" | "fono" | "}"
# (i.e. one opening/closing brace is dropped and "Tel" is re-split so the
# leading "T" becomes lower-case "t" while "el" stays together).
#
# A plain Find/Replace (or any Range.Text assignment / Delete / InsertText)
# collapses every run it touches into a single run, so it cannot reproduce
# that run layout by itself. The trick used below is:
#   1) Locate the placeholder and rewrite its text content in one shot
#      (this naturally merges it into a single run).
#   2) Re-impose the required run boundaries by toggling a character
#      formatting property (Bold) on and back off over each exact final
#      segment. Word keeps a distinct run for a sub-range whose formatting
#      was explicitly touched, even though the end result is visually
#      identical to its neighbours, which lets us recreate the six runs
#      from the diff without changing anything else about the text.
# ------------------------------------------------------------------------

$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute("{{Telefono}}", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)

if ($found) {
    $start = $findRange.Start
    $end = $findRange.End

    # Step 1: replace the whole placeholder text (merges into one run).
    $whole = $d.Range($start, $end)
    $whole.Text = "{telefono}"

    # Step 2: force the six target run boundaries back into existence:
    #   "{"  "t"  "el"  "e"  "fono"  "}"
    $segmentLengths = @(1, 1, 2, 1, 4, 1)
    $pos = $start
    foreach ($len in $segmentLengths) {
        $segment = $d.Range($pos, $pos + $len)
        $segment.Font.Bold = 1
        $segment.Font.Bold = 0
        $pos = $pos + $len
    }
}
